# Remove the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph, the
# "(c) 2020 . Contact: ..." paragraph, and the blank paragraph that sits
# between them and the "LOM3013: ..." requirement line, while leaving the
# trailing blank paragraph (just before the page-break paragraph) intact.

$d = $word.ActiveDocument

$startPara = $null
$endPara = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $startPara = $p
    }
    if ($t -like "*Contact: luizeleno@usp.br*") {
        $endPara = $p
    }
}

# Include the blank paragraph immediately preceding the "Ver no Jupiter..."
# paragraph so all three paragraphs (and their paragraph marks) are removed.
$deleteStart = $startPara.Previous().Range.Start
$deleteEnd = $endPara.Range.End

$r = $d.Range($deleteStart, $deleteEnd)
$r.Delete()
